$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Raw Data 2" table (A32:B41) got new measurements and grew by one row
# (A33:B42). Update the existing rows' values first.
$data = @(
    @(7.5, 1),
    @(9.5, 1),
    @(10, 4),
    @(10.5, 1),
    @(11, 3),
    @(11.5, 5),
    @(12.5, 1),
    @(13, 1),
    @(14, 2)
)

$row = 33
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Add the new trailing row (42), matching the formatting of the row above it.
$ws.Range("A41:B41").Copy()
$ws.Range("A42:B42").PasteSpecial(-4122)
$ws.Cells.Item(42, 1).Value = 18
$ws.Cells.Item(42, 2).Value = 1
$ws.Rows.Item(42).RowHeight = 15

# Leave the selection where the edit ended up, scrolled down to show it.
$ws.Range("B42").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 26
